$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = 45763
$ws.Range("B6").Value = 5.2
$ws.Range("C6").Value = 5.2
$ws.Range("D6").Value = 5.2
$ws.Range("E6").Value = 5.3
$ws.Range("F6").Value = 5.7
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 6.7
$ws.Range("I6").Value = 5.5
$ws.Range("J6").Value = 4.7
$ws.Range("K6").Value = 11.6
$ws.Range("L6").Value = 8.2
$ws.Range("M6").Value = 6.5
